# edit.ps1
# Applies the "Better AI, Robot's Layers" changes:
#  1. Merge the two runs "6" + " = Leg" into a single run "6 = Leg".
#  2. Merge "Body"+"1" / "Body"+"2" / "Body"+"3" run-pairs into single runs.
#  3. Append a new "Слои роботов" section (two blank paragraphs, a heading
#     paragraph, and a new 5x5 table) right before the trailing bookmarked
#     paragraph, then drop the now-superfluous "_GoBack" bookmark so the
#     document ends with a single empty paragraph (the bookmark now lives
#     inside the new table, on the "Оружие П" cell of the last row).

$d = $word.ActiveDocument

# --- 1) "6" + " = Leg" -> "6 = Leg" -------------------------------------
$rng = $d.Content
$rng.Find.Execute("6 = Leg", $true, $false, $false, $false, $false, $true, 1, $false, "6 = Leg", 2) | Out-Null

# --- 2) "Body" + "1"/"2"/"3" -> "Body1"/"Body2"/"Body3" -----------------
foreach ($t in @("Body1", "Body2", "Body3")) {
    $rng = $d.Content
    $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2) | Out-Null
}

# --- 3) New "Слои роботов" section + table ------------------------------
# Insert right at the start of the document's current last paragraph (the
# one holding the lone "_GoBack" bookmark). The new table below carries its
# own "_GoBack" bookmark (on the "Оружие П" cell of the last row), so the
# old bookmark -- now a duplicate -- is dropped automatically, leaving the
# old paragraph as a plain, empty paragraph. A temporary placeholder run is
# used as the very last paragraph of the inserted fragment so that the
# (otherwise completely empty) trailing paragraph is not optimized away
# during the insert; it is emptied back out right afterwards.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$insertPoint = $lastPara.Range.Duplicate
$insertPoint.Collapse(1)

$newSectionXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml">
<pkg:xmlData>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Слои роботов</w:t>
  </w:r>
</w:p>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="a3"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="1689"/>
    <w:gridCol w:w="1914"/>
    <w:gridCol w:w="1914"/>
    <w:gridCol w:w="1914"/>
    <w:gridCol w:w="1914"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1689" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Вверх</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Вниз</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Влево</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Вправо</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1689" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>Top 5</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Корпус</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Оружие 2</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Оружие Л</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Оружие П</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1689" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>MidTop</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve"> 4</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Ноги</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Корпус</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Ноги</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Ноги</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1689" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>MidBot</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve"> 3</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Оружие 2</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Ноги</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Корпус</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Корпус</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1689" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="C5E0B3" w:themeFill="accent6" w:themeFillTint="66"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:val="en-US"/>
          </w:rPr>
          <w:t>Bot 2</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p/>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Оружие П</w:t>
        </w:r>
        <w:bookmarkStart w:id="0" w:name="_GoBack"/>
        <w:bookmarkEnd w:id="0"/>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1914" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Оружие Л</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>ZZ_TRAILING_PLACEHOLDER_ZZ</w:t>
  </w:r>
</w:p>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertPoint.InsertXML($newSectionXml) | Out-Null

$placeholder = $d.Content
$placeholder.Find.Execute("ZZ_TRAILING_PLACEHOLDER_ZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
